$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# --- Update hotel_info: insert a new "State" column between Hotel_Name and City ---
$hotel = $wb.Worksheets.Item("hotel_info")
$hotel.Columns.Item(3).Insert()
$hotel.Range("C1").Value = "State"
$hotel.Range("C2").Value = "Louisiana"

# --- Reorder worksheet tabs: review_info first, hotel_info second ---
$review = $wb.Worksheets.Item("review_info")
$review.Move($wb.Worksheets.Item(1))

# Keep the first tab active/selected
$wb.Worksheets.Item(1).Activate()
